$d = $word.ActiveDocument

# 1. Update the title and intro paragraph text.
[void]$d.Content.Find.Execute("Complex Test Document", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "Test Document with Table", 2)
[void]$d.Content.Find.Execute("This document has multiple tables.", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "This is a test document.", 2)

# 2. Remove the first table (the Key/Value table) entirely.
$d.Tables.Item(1).Delete()

# 3. Remove the now-orphaned "Here is another table:" paragraph that
#    used to introduce the second table.
foreach ($p in $d.Content.Paragraphs) {
    if ($p.Range.Text -like "Here is another table:*") {
        $p.Range.Delete()
        break
    }
}

# 4. Apply the "Light Grid Accent 1" table style to the remaining table.
$t = $d.Tables.Item(1)
$t.Style = "Light Grid Accent 1"

# 5. Rename the remaining table's headers and data
#    (Product/Price/Stock + Apple/Banana rows -> Name/Age/City + Alice/Bob rows).
[void]$d.Content.Find.Execute("Product", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "Name", 2)
[void]$d.Content.Find.Execute("Price", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "Age", 2)
[void]$d.Content.Find.Execute("Stock", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "City", 2)
[void]$d.Content.Find.Execute("Apple", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "Alice", 2)
[void]$d.Content.Find.Execute("1.99", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "30", 2)
[void]$d.Content.Find.Execute("100", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "NYC", 2)
[void]$d.Content.Find.Execute("Banana", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "Bob", 2)
[void]$d.Content.Find.Execute("0.99", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "25", 2)
[void]$d.Content.Find.Execute("50", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "LA", 2)

# 6. Drop the last data row (the "Orange" row) so only Alice/Bob remain.
$t = $d.Tables.Item(1)
$t.Rows.Item($t.Rows.Count).Delete()
